$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("standard_curve_no2-test")

# --- Turn the existing one-off formulas into filled/shared formula ranges ---
$ws.Range("F2:F12").Formula = "=64.069*AVERAGE(D2:D2)-2.8707"
$ws.Range("J2:J8").Formula = "=SUM(H2:I2)/3"

# F13 keeps its own (non-shared) formula, now multiplying by E13 instead of a
# hard-coded 4
$ws.Range("F13").Formula = "=(64.069*AVERAGE(D13:D13)-2.8707)*E13"

# --- New sample rows 14-37 (first no2- data) ---
$ws.Range("C14").Value = "P1T0"
$ws.Range("D14").Value = 0.126
$ws.Range("E14").Value = 4

$ws.Range("C15").Value = "P1T3"
$ws.Range("D15").Value = 0.198
$ws.Range("E15").Value = 4

$ws.Range("C16").Value = "P1T8"
$ws.Range("D16").Value = 0.19
$ws.Range("E16").Value = 4

$ws.Range("C17").Value = "P1T15"
$ws.Range("D17").Value = 0.179
$ws.Range("E17").Value = 4

$ws.Range("C18").Value = "P1T23"
$ws.Range("D18").Value = 0.197
$ws.Range("E18").Value = 4

$ws.Range("C19").Value = "P1T27"
$ws.Range("D19").Value = 0.183
$ws.Range("E19").Value = 4

$ws.Range("C20").Value = "P2T0"
$ws.Range("D20").Value = 0.035
$ws.Range("E20").Value = 4

$ws.Range("C21").Value = "P2T3"
$ws.Range("D21").Value = 0.814
$ws.Range("E21").Value = 20

$ws.Range("C22").Value = "P2T8"
$ws.Range("D22").Value = 0.683
$ws.Range("E22").Value = 20

$ws.Range("C23").Value = "P2T15"
$ws.Range("D23").Value = 0.51
$ws.Range("E23").Value = 10

$ws.Range("C24").Value = "P2T23"
$ws.Range("D24").Value = 0.382
$ws.Range("E24").Value = 10

$ws.Range("C25").Value = "P2T26"
$ws.Range("D25").Value = 0.31
$ws.Range("E25").Value = 10

$ws.Range("C26").Value = "P3T0"
$ws.Range("D26").Value = 0.042
$ws.Range("E26").Value = 4

$ws.Range("C27").Value = "P3T3"
$ws.Range("D27").Value = 0.679
$ws.Range("E27").Value = 20

$ws.Range("C28").Value = "P3T8"
$ws.Range("D28").Value = 0.629
$ws.Range("E28").Value = 20

$ws.Range("C29").Value = "P3T15"
$ws.Range("D29").Value = 0.807
$ws.Range("E29").Value = 10

$ws.Range("C30").Value = "P3T23"
$ws.Range("D30").Value = 0.626
$ws.Range("E30").Value = 10

$ws.Range("C31").Value = "P3T27"
$ws.Range("D31").Value = 0.558
$ws.Range("E31").Value = 10

$ws.Range("C32").Value = "P4T0"
$ws.Range("D32").Value = 0.055
$ws.Range("E32").Value = 1

$ws.Range("C33").Value = "P4T3"
$ws.Range("D33").Value = 0.036
$ws.Range("E33").Value = 1

$ws.Range("C34").Value = "P4T8"
$ws.Range("D34").Value = 0.035
$ws.Range("E34").Value = 1

$ws.Range("C35").Value = "P4T15"
$ws.Range("D35").Value = 0.037
$ws.Range("E35").Value = 1

$ws.Range("C36").Value = "P4T23"
$ws.Range("D36").Value = 0.038
$ws.Range("E36").Value = 1

$ws.Range("C37").Value = "P4T27"
$ws.Range("D37").Value = 0.034
$ws.Range("E37").Value = 1

# Fill formula down column F for the new rows (one shared formula group)
$ws.Range("F14:F37").Formula = "=(64.069*AVERAGE(D14:D14)-2.8707)*E14"

# The last four "plates" (P4Txx) samples get a black (explicit RGB) font
# color instead of the theme-based default
$ws.Range("C32:C37").Font.Color = 0

# Re-apply the "F"+batch+"-"+date formula over C3:C13 (narrowed down from the
# previous C3:C17 range now that C14:C37 hold plain sample names)
$ws.Range("C3:C13").Formula = '=_xlfn.CONCAT("F",B3,"-",A3)'

# --- Selection / view bookkeeping ---
$ws.Range("E23").Select()
